$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.612.44'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.507.69'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.56'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.33'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.506.41'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.94%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.26'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.106.84'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.66'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.509.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.602.54'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.93'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '391.39'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.650.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.93'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.57'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.45'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.41%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.513.27'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '24.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '170.97'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.96'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.813'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.42'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.55%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.08'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.40'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.453.80'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.895'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.25%  '
